$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 311, shifting rows 311:401 down to 312:402
$ws.Rows.Item(311).Insert()

# Populate the new row 311 with the new record
# Columns A,B,C,E,F,G,H,I,J,Q are constant across the dataset; copy from row 312 (the row pushed down)
$ws.Cells.Item(311,1).Value2  = $ws.Cells.Item(312,1).Value2   # A Mercado ID
$ws.Cells.Item(311,2).Value2  = $ws.Cells.Item(312,2).Value2   # B Mercado
$ws.Cells.Item(311,3).Value2  = $ws.Cells.Item(312,3).Value2   # C Region
$ws.Cells.Item(311,4).Value2  = 44559                          # D Fecha
$ws.Cells.Item(311,5).Value2  = $ws.Cells.Item(312,5).Value2   # E Codreg
$ws.Cells.Item(311,6).Value2  = $ws.Cells.Item(312,6).Value2   # F Tipo
$ws.Cells.Item(311,7).Value2  = $ws.Cells.Item(312,7).Value2   # G Producto ID
$ws.Cells.Item(311,8).Value2  = $ws.Cells.Item(312,8).Value2   # H Producto
$ws.Cells.Item(311,9).Value2  = $ws.Cells.Item(312,9).Value2   # I Categoria ID
$ws.Cells.Item(311,10).Value2 = $ws.Cells.Item(312,10).Value2  # J Categoria
$ws.Cells.Item(311,11).Value2 = "Valencia"                     # K Variedad
$ws.Cells.Item(311,12).Value2 = "Primera"                      # L Calidad
$ws.Cells.Item(311,13).Value2 = 250                             # M Volumen
$ws.Cells.Item(311,14).Value2 = 8000                            # N Precio minimo
$ws.Cells.Item(311,15).Value2 = 8000                            # O Precio maximo
$ws.Cells.Item(311,16).Value2 = 8000                            # P Precio promedio ponderado
$ws.Cells.Item(311,17).Value2 = $ws.Cells.Item(312,17).Value2  # Q Unidad de comercializacion
$ws.Cells.Item(311,18).Value2 = "Provincia de Melipilla"       # R Origen
$ws.Cells.Item(311,19).Value2 = 533                             # S Precio $/Kg
$ws.Cells.Item(311,20).Value2 = 15                              # T Kg / unidad
